$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "53.173.47"
$ws.Range("E2").Value = "  -12.51%  "

# Row 3
$ws.Range("D3").Value = "2.327.10"
$ws.Range("E3").Value = "  -19.91%  "

# Row 4
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "438.02"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -17.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "121.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -15.49%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.474"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -14.82%  "

# Row 9
$ws.Range("D9").Value = "2.331.44"
$ws.Range("E9").Value = "  -20.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0915"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -15.47%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.27"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -12.99%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.309"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -14.58%  "

# Row 13
$ws.Range("E13").Value = "  -3.80%  "

# Row 14
$ws.Range("D14").Value = "2.727.96"
$ws.Range("E14").Value = "  -20.12%  "

# Row 15
$ws.Range("D15").Value = "53.164.00"
$ws.Range("E15").Value = "  -12.44%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -16.28%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000120"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -15.41%  "

# Row 18
$ws.Range("D18").Value = "2.340.10"
$ws.Range("E18").Value = "  -19.59%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -21.15%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "301.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -16.91%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -21.97%  "

# Row 22
$ws.Range("E22").Value = "  -0.24%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.15%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -18.85%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "55.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -14.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.153"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -15.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.369"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -19.00%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -11.41%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0704"
$ws.Range("E31").Value = "  -17.94%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "144.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.21%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "17.21"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -13.14%  "

# Row 34
$ws.Range("E34").Value = "  -19.87%  "

# Row 35
$ws.Range("E35").Value = "  -15.34%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.55"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -18.92%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.832"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -17.23%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -17.05%  "

# Row 39
$ws.Range("B39").Value = "OKB"
$ws.Range("C39").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "33.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.97%  "

# Row 40
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.994"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.28"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.29%  "

# Row 42
$ws.Range("E42").Value = "  -14.81%  "

# Row 43
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0501"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -14.67%  "

# Row 44
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -18.95%  "

# Row 45
$ws.Range("D45").Value = "1.908.28"
$ws.Range("E45").Value = "  -17.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.524"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -19.45%  "

# Row 47
$ws.Range("E47").Value = "  -11.76%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0832"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -10.16%  "

# Row 49
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "15.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -23.68%  "

# Row 50
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.91"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -21.74%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.40%  "
